# Commit: "inizio modifiche per acquisire segnali da 2 macchine"
# Clears the header-row labels (A1:FG1) and refreshes several sampled
# sensor readings in row 2 now that data is being merged from 2 machines.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: wipe out the old single-machine column headers (A1:FG1). ---
# The remaining headers from FH1 onward were already blank and stay that way.
$ws.Range("A1:FG1").ClearContents()

# --- Row 2: updated sampled values reflecting the 2-machine acquisition. ---
$ws.Range("G2").Value = 2162733
$ws.Range("H2").Value = 2162733

$ws.Range("Q2").Value = 45
$ws.Range("R2").Value = 25
$ws.Range("S2").Value = 34

$ws.Range("V2").Value = 31
$ws.Range("W2").Value = 40
$ws.Range("X2").Value = 32
$ws.Range("Y2").Value = 39
$ws.Range("Z2").Value = 26

$ws.Range("AA2").Value = 1
$ws.Range("AB2").Value = 2
$ws.Range("AC2").Value = 230
$ws.Range("AD2").Value = 3
$ws.Range("AE2").Value = 483
$ws.Range("AF2").Value = 0
$ws.Range("AH2").Value = 0

$ws.Range("DU2").Value = "O4084"
$ws.Range("DV2").Value = "O5084"

# These columns hold purely numeric-looking text (inline strings) in the
# sheet, so force the number format to text before assigning -- otherwise
# Excel would coerce the numeric-looking string back into a real number.
$ws.Range("EB2").NumberFormat = "@"
$ws.Range("EB2").Value = "2"
$ws.Range("EC2").NumberFormat = "@"
$ws.Range("EC2").Value = "25"
$ws.Range("ED2").NumberFormat = "@"
$ws.Range("ED2").Value = "34"

$ws.Range("EJ2").NumberFormat = "@"
$ws.Range("EJ2").Value = "26"
$ws.Range("EK2").NumberFormat = "@"
$ws.Range("EK2").Value = "250"

$ws.Range("FA2").Value = 27
$ws.Range("FB2").Value = 29
$ws.Range("FC2").Value = 27
$ws.Range("FD2").Value = 29
$ws.Range("FE2").Value = 2539
$ws.Range("FF2").Value = 6
$ws.Range("FG2").Value = 17
